$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "66.782.48"
Set-TextValue $ws.Range("E2") "  -1.28%  "
Set-TextValue $ws.Range("D3") "3.495.93"
Set-TextValue $ws.Range("E3") "  -0.52%  "
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "602.16"
Set-TextValue $ws.Range("E5") "  -1.38%  "
Set-TextValue $ws.Range("D6") "147.71"
Set-TextValue $ws.Range("E6") "  -3.12%  "
Set-TextValue $ws.Range("D7") "3.495.36"
Set-TextValue $ws.Range("E7") "  -0.51%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("E9") "  -1.75%  "
Set-TextValue $ws.Range("E10") "  -1.18%  "
Set-TextValue $ws.Range("D11") "7.88"
Set-TextValue $ws.Range("E11") "  +2.96%  "
Set-TextValue $ws.Range("E12") "  -2.41%  "
Set-TextValue $ws.Range("E13") "  -1.98%  "
Set-TextValue $ws.Range("D14") "4.085.92"
Set-TextValue $ws.Range("E14") "  -0.49%  "
Set-TextValue $ws.Range("D15") "31.25"
Set-TextValue $ws.Range("E15") "  -4.55%  "
Set-TextValue $ws.Range("D16") "3.490.78"
Set-TextValue $ws.Range("E16") "  -0.58%  "
Set-TextValue $ws.Range("D17") "66.762.55"
Set-TextValue $ws.Range("E17") "  -1.24%  "
Set-TextValue $ws.Range("E18") "  -0.62%  "
Set-TextValue $ws.Range("D19") "10.49"
Set-TextValue $ws.Range("E19") "  +6.19%  "
Set-TextValue $ws.Range("D20") "6.37"
Set-TextValue $ws.Range("E20") "  -3.12%  "
Set-TextValue $ws.Range("D21") "15.33"
Set-TextValue $ws.Range("E21") "  -1.96%  "
Set-TextValue $ws.Range("D22") "433.68"
Set-TextValue $ws.Range("E22") "  -3.62%  "
Set-TextValue $ws.Range("D23") "0.607"
Set-TextValue $ws.Range("E23") "  -4.36%  "
Set-TextValue $ws.Range("D24") "79.74"
Set-TextValue $ws.Range("E24") "  +1.94%  "
Set-TextValue $ws.Range("D25") "3.632.23"
Set-TextValue $ws.Range("E25") "  -0.53%  "
Set-TextValue $ws.Range("E26") "  -0.11%  "
Set-TextValue $ws.Range("E27") "  -3.15%  "
Set-TextValue $ws.Range("E28") "  -7.39%  "
Set-TextValue $ws.Range("D29") "9.82"
Set-TextValue $ws.Range("E29") "  -2.93%  "
Set-TextValue $ws.Range("D30") "8.22"
Set-TextValue $ws.Range("E30") "  -7.90%  "
Set-TextValue $ws.Range("E31") "  -1.17%  "
Set-TextValue $ws.Range("D32") "1.59"
Set-TextValue $ws.Range("E32") "  -4.05%  "
Set-TextValue $ws.Range("D33") "1.00"
Set-TextValue $ws.Range("E33") "  +0.06%  "
Set-TextValue $ws.Range("D34") "0.165"
Set-TextValue $ws.Range("E34") "  -2.36%  "
Set-TextValue $ws.Range("D35") "25.35"
Set-TextValue $ws.Range("E35") "  -1.77%  "
Set-TextValue $ws.Range("D36") "3.487.75"
Set-TextValue $ws.Range("E36") "  -0.53%  "
Set-TextValue $ws.Range("B37") "ImmutableX"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "1.79"
Set-TextValue $ws.Range("E37") "  -4.86%  "
Set-TextValue $ws.Range("B38") "NEARProtocol"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D38") "5.87"
Set-TextValue $ws.Range("E38") "  -5.47%  "
Set-TextValue $ws.Range("D39") "7.99"
Set-TextValue $ws.Range("E39") "  -1.15%  "
Set-TextValue $ws.Range("E41") "  -0.06%  "
Set-TextValue $ws.Range("D42") "0.0890"
Set-TextValue $ws.Range("E42") "  -1.11%  "
Set-TextValue $ws.Range("D43") "170.25"
Set-TextValue $ws.Range("E43") "  -2.37%  "
Set-TextValue $ws.Range("E44") "  -9.17%  "
Set-TextValue $ws.Range("D45") "5.41"
Set-TextValue $ws.Range("E45") "  -1.61%  "
Set-TextValue $ws.Range("D46") "0.897"
Set-TextValue $ws.Range("E46") "  +1.47%  "
Set-TextValue $ws.Range("B47") "InjectiveProtocol"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D47") "28.75"
Set-TextValue $ws.Range("E47") "  -5.02%  "
Set-TextValue $ws.Range("B48") "OKB"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D48") "45.69"
Set-TextValue $ws.Range("E48") "  -2.28%  "
Set-TextValue $ws.Range("D49") "1.32"
Set-TextValue $ws.Range("E49") "  +0.67%  "
Set-TextValue $ws.Range("E50") "  -3.13%  "
Set-TextValue $ws.Range("E51") "  -4.77%  "
